# Mise a jour 2024
# - Bump the "last updated" date field (28/01/2023 -> 14/01/2024) everywhere
#   it appears: the slide master, every slide layout and the notes master.
# - Bump the Android/iOS version call-outs on slide 2.

$p = $ppt.ActivePresentation

$oldDate = "28/01/2023"
$newDate = "14/01/2024"

function Update-DateField {
    param($shapes)

    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master date placeholder.
$master = $p.SlideMaster
Update-DateField $master.Shapes

# Every custom layout's date placeholder.
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Update-DateField $layout.Shapes
}

# Notes master date placeholder.
$notesMaster = $p.NotesMaster
Update-DateField $notesMaster.Shapes

# Slide 2: bump the Android/iOS version mentions.
$slide2 = $p.Slides.Item(2)
$contentShape = $slide2.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$fullText = $tr.Text

$oldAndroid = "Version 13 en 2022"
$newAndroid = "Version 14 en 2024"
$idxAndroid = $fullText.IndexOf($oldAndroid)
if ($idxAndroid -ge 0) {
    $chars = $tr.Characters($idxAndroid + 1, $oldAndroid.Length)
    $chars.Text = $newAndroid
}

$fullText = $tr.Text
$oldIos = "Version 16 en 2023"
$newIos = "Version 17 en 2023"
$idxIos = $fullText.IndexOf($oldIos)
if ($idxIos -ge 0) {
    $chars = $tr.Characters($idxIos + 1, $oldIos.Length)
    $chars.Text = $newIos
}
